$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The original column A held a running row index (1,2,3,4) that is no longer
# needed; removing it shifts the real table (B:D -> A:C).
$ws.Columns.Item(1).Delete()

# Insert three new rows, one above each existing time-slot row, to hold the
# newly-recorded 10-minute slots (16.00-16.10, 16.20-16.30, 16.40-16.50).
# This also pushes the 4 trailing blank rows down from 6-9 to 9-12.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(7).Insert()

# Re-write the header rows (unchanged content, now in columns A:C).
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"

$ws.Range("A2").Value = "JAM"
$ws.Range("B2").Value = "Jumlah Motor"
$ws.Range("C2").Value = "Jumlah Mobil"

# Full, chronologically-ordered data table.
$ws.Range("A3").Value = "16.00 - 16.10"
$ws.Range("B3").Value = 169
$ws.Range("C3").Value = 4

$ws.Range("A4").Value = "16.10 - 16.20"
$ws.Range("B4").Value = 128
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "16.20 - 16.30"
$ws.Range("B5").Value = 156
$ws.Range("C5").Value = 1

$ws.Range("A6").Value = "16.30 - 16.40"
$ws.Range("B6").Value = 151
$ws.Range("C6").Value = 3

$ws.Range("A7").Value = "16.40 - 16.50"
$ws.Range("B7").Value = 141
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = "16.50 - 17.00"
$ws.Range("B8").Value = 134
$ws.Range("C8").Value = 4

# Match the saved selection / active cell from the edited workbook.
$ws.Range("C7").Select()
